$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Mon Oct 09 22:46:38 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 22:46:52 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 22:47:06 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 22:47:19 EDT 2023"
